$d = $word.ActiveDocument

# --- 1) Collapse "Resume at " + word-ladder hyperlink + blank paragraph
#        into the "come back to:" paragraph, which becomes
#        "Main List is done. Come back to:" (two runs, matching the diff).

$start = $d.Paragraphs.Item(6).Range.Start
$end = $d.Paragraphs.Item(9).Range.Start
$rng = $d.Range($start, $end)
$rng.Delete()

$p = $d.Paragraphs.Item(6)
$r = $p.Range
# Replace only the paragraph's text content, leaving its paragraph mark intact
$textRange = $d.Range($r.Start, $r.End - 1)
$xml = '<?xml version="1.0"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Main List is done. C</w:t></w:r><w:r><w:t>ome back to:</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$textRange.InsertXML($xml)

# --- 2) Drop two of the three trailing empty paragraphs at the end of the
#        document, leaving a single blank paragraph before the sectPr.

for ($k = 0; $k -lt 2; $k++) {
    $n = $d.Paragraphs.Count
    $p = $d.Paragraphs.Item($n - 1)
    $rng = $d.Range($p.Range.Start, $p.Range.Start + 1)
    $rng.Delete()
}
